$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CasesTab query (B2): remove TRIM() wraps on Case ID/Study Code/Study Type,
# and change the Age column to a text-preserving CASE expression.
$ws.Range("B2").Value2 = "SELECT`r`n    DISTINCT c.case_record_id AS `"Case ID`",`r`n    st.clinical_study_designation AS `"Study Code`",`r`n    st.clinical_study_type AS `"Study Type`",`r`n    dmg.breed AS `"Breed`",`r`n    diag.disease_term AS `"Diagnosis`",`r`n    diag.stage_of_disease AS `"Stage Of Disease`",`r`n    CASE `r`n    WHEN dmg.patient_age_at_enrollment = CAST(dmg.patient_age_at_enrollment AS INT) `r`n    THEN CAST(CAST(dmg.patient_age_at_enrollment AS INT) AS TEXT)`r`n    ELSE CAST(dmg.patient_age_at_enrollment AS TEXT)`r`nEND AS `"Age`", `r`n    COALESCE(TRIM(dmg.sex), '') AS `"Sex`",   `r`n    COALESCE(dmg.neutered_indicator, '') AS `"Neutered Status`",  `r`n    COALESCE(`r`n        CASE `r`n            WHEN dmg.weight = CAST(dmg.weight AS INT) THEN CAST(dmg.weight AS INT)`r`n            ELSE dmg.weight `r`n        END, `r`n    '') AS `"Weight (kg)`",        `r`n    COALESCE(diag.best_response, '') AS `"Response to Treatment`",  `r`n    COALESCE(coh.cohort_description, '') AS `"Cohort`"  `r`nFROM `r`n    df_case c`r`nJOIN `r`n    df_study st ON c.`"study.clinical_study_designation`" = st.clinical_study_designation`r`nJOIN `r`n    df_demographic dmg ON dmg.`"case.case_record_id`" = c.case_record_id`r`nJOIN `r`n    df_diagnosis diag ON diag.`"case.case_record_id`" = c.case_record_id`r`nJOIN `r`n    df_enrollment enr ON enr.`"case.case_record_id`" = c.case_record_id`r`nJOIN `r`n    df_program p ON st.`"program.program_acronym`" = p.program_acronym`r`nJOIN `r`n    df_sample smp ON smp.`"case.case_record_id`" = c.case_record_id`r`nJOIN `r`n    df_publication pub ON pub.`"study.clinical_study_designation`" = st.clinical_study_designation`r`nLEFT JOIN `r`n    df_case_file cf ON cf.`"sample.sample_id`" = smp.sample_id`r`nLEFT JOIN `r`n    df_study_file sf ON sf.`"study.clinical_study_designation`" = st.clinical_study_designation`r`nLEFT JOIN`r`n    df_cohort coh ON coh.`"study.clinical_study_designation`" = st.clinical_study_designation`r`nWHERE`r`n    st.clinical_study_designation = 'TCL01' AND dmg.breed = 'Golden Retriever'`r`nORDER BY `r`n    c.case_record_id ASC`r`nLIMIT 100;"

# Update SamplesTab query (B3): remove TRIM() wrap on necropsy_sample.
$ws.Range("B3").Value2 = "SELECT DISTINCT`r`n    smp.sample_id AS `"Sample ID`",`r`n    c.case_record_id AS `"Case ID`",`r`n    COALESCE(dmg.breed, '') AS `"Breed`",`r`n    COALESCE(diag.disease_term, '') AS `"Diagnosis`",`r`n    COALESCE(smp.sample_site, '') AS `"Sample Site`",`r`n    COALESCE(smp.summarized_sample_type, '') AS `"Sample Type`",`r`n    COALESCE(smp.specific_sample_pathology, '') AS `"Pathology/Morphology`",`r`n    COALESCE(smp.tumor_grade, '') AS `"Tumor Grade`",`r`n    COALESCE(smp.sample_chronology, '') AS `"Sample Chronology`",`r`n    COALESCE(smp.percentage_tumor, '') AS `"Percentage Tumor`",`r`n    COALESCE(smp.necropsy_sample, '') AS `"Necropsy Sample`",`r`n    COALESCE(smp.sample_preservation, '') AS `"Sample Preservation`"`r`nFROM `r`n    df_sample smp`r`nJOIN `r`n    df_case c ON smp.`"case.case_record_id`" = c.case_record_id`r`nJOIN `r`n    df_publication pub ON pub.`"study.clinical_study_designation`" = st.clinical_study_designation`r`nJOIN `r`n    df_demographic dmg ON dmg.`"case.case_record_id`" = c.case_record_id`r`nJOIN `r`n    df_diagnosis diag ON diag.`"case.case_record_id`" = c.case_record_id`r`nJOIN `r`n    df_enrollment enr ON enr.`"case.case_record_id`" = c.case_record_id`r`nJOIN `r`n    df_program p ON st.`"program.program_acronym`" = p.program_acronym`r`nJOIN `r`n    df_study st ON c.`"study.clinical_study_designation`" = st.clinical_study_designation`r`nLEFT JOIN `r`n    df_case_file cf ON cf.`"sample.sample_id`" = smp.sample_id`r`nLEFT JOIN `r`n    df_study_file sf ON sf.`"study.clinical_study_designation`" = st.clinical_study_designation`r`nWHERE `r`n   st.clinical_study_designation = 'TCL01' AND dmg.breed = 'Golden Retriever'`r`nORDER BY `r`n    smp.sample_id ASC`r`nLIMIT 100;"

# Update StudyFilesTab query (B5): remove TRIM() wraps on file_name/file_type/file_description.
$ws.Range("B5").Value2 = "SELECT DISTINCT`r`n    sf.file_name AS `"File Name`",`r`n    sf.file_type AS `"File Type`",`r`n    'study' AS `"Association`",`r`n    sf.file_description AS `"Description`",`r`n    CASE`r`n        WHEN sf.file_name LIKE '%.bai' THEN 'bai'`r`n        WHEN sf.file_name LIKE '%.bam' THEN 'bam'`r`n        WHEN sf.file_name LIKE '%.csv' THEN 'csv'`r`n        WHEN sf.file_name LIKE '%.doc' THEN 'doc'`r`n        WHEN sf.file_name LIKE '%.docx' THEN 'docx'`r`n        WHEN sf.file_name LIKE '%.gz' THEN 'gz'`r`n        WHEN sf.file_name LIKE '%.pdf' THEN 'pdf'`r`n        WHEN sf.file_name LIKE '%.rtf' THEN 'rtf'`r`n        WHEN sf.file_name LIKE '%.tbi' THEN 'tbi'`r`n        WHEN sf.file_name LIKE '%.tif' THEN 'tif'`r`n        WHEN sf.file_name LIKE '%.xls' THEN 'xls'`r`n        WHEN sf.file_name LIKE '%.xlsx' THEN 'xlsx'`r`n        ELSE 'Unknown'`r`n    END AS `"Format`",`r`n    CASE     `r`n        WHEN sf.file_size >= 1024 * 1024 * 1024 THEN `r`n            ROUND(sf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB' `r`n        WHEN sf.file_size >= 1024 * 1024 THEN `r`n            ROUND(sf.file_size / (1024.0 * 1024.0), 2) || ' MB' `r`n        WHEN sf.file_size >= 1024 THEN `r`n            ROUND(sf.file_size / 1024.0, 2) || ' KB' `r`n        ELSE `r`n            ROUND(sf.file_size, 2) || ' Bytes' `r`n    END AS `"Size`",`r`n    st.clinical_study_designation AS `"Study Code`"`r`nFROM `r`n    df_case_file cf`r`nJOIN `r`n    df_sample smp ON cf.`"sample.sample_id`" = smp.sample_id`r`nJOIN `r`n    df_case c ON smp.`"case.case_record_id`" = c.case_record_id`r`nJOIN `r`n    df_study st ON c.`"study.clinical_study_designation`" = st.clinical_study_designation`r`nJOIN `r`n    df_program p ON st.`"program.program_acronym`" = p.program_acronym`r`nJOIN `r`n    df_demographic dmg ON dmg.`"case.case_record_id`" = c.case_record_id`r`nJOIN `r`n    df_diagnosis diag ON diag.`"case.case_record_id`" = c.case_record_id`r`nJOIN `r`n    df_enrollment enr ON enr.`"case.case_record_id`" = c.case_record_id`r`nJOIN `r`n    df_publication pub ON pub.`"study.clinical_study_designation`" = st.clinical_study_designation`r`nLEFT JOIN `r`n    df_study_file sf ON sf.`"study.clinical_study_designation`" = st.clinical_study_designation`r`nWHERE`r`n    st.clinical_study_designation = 'TCL01' AND dmg.breed = 'Golden Retriever'`r`nORDER BY `r`n    sf.file_name ASC`r`nLIMIT 100;"

# Move the active selection/scroll position to B5 (StudyFilesTab row), matching the author's
# final cursor position after editing the StudyFilesTab query last.
$ws.Range("A5").Select()
$ws.Range("B5").Select()
